$d = $word.ActiveDocument
$bullet = [char]0x2022

# The "Programa" section has two paragraphs (PT and EN) whose text runs
# together several bullet ("•") separated topics inside a single <w:t>.
# Split each topic onto its own line by inserting a manual line break
# (w:br, "^l") right before every bullet except the very first one in the
# paragraph. The wildcard pattern "(?)" captures exactly one preceding
# character, so it cannot match at the very start of the paragraph (where
# the bullet has no predecessor) - this naturally skips the leading
# bullet while catching all the others.
#
# The "Bibliografia" paragraph lists four numbered references glued
# together (e.g. "...2005.2.LIMA..."). Insert a line break before each of
# the three later entry markers (2., 3., 4.) while leaving the first
# entry intact.
#
# Paragraphs are located by their known leading text rather than a fixed
# index, so the script is robust even if paragraph numbering shifts.

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i).Range
    $text = $para.Text

    if ($text.StartsWith($bullet)) {
        $para.Find.ClearFormatting()
        $para.Find.Replacement.ClearFormatting()
        $ok = $para.Find.Execute("(?)" + $bullet, $true, $false, $true, $false, $false, $true, 1, $false, "\1^l" + $bullet, 2)
        if (-not $ok) {
            throw "Bullet split failed for paragraph $i"
        }
    }
    elseif ($text.StartsWith("1.CAMARGO")) {
        $para.Find.ClearFormatting()
        $para.Find.Replacement.ClearFormatting()
        $ok = $para.Find.Execute("(.)([2-9]\.)", $true, $false, $true, $false, $false, $true, 1, $false, "\1^l\2", 2)
        if (-not $ok) {
            throw "Reference split failed for Bibliografia paragraph $i"
        }
    }
}
